$wb = $excel.ActiveWorkbook

# ---- Add the four new sheets in order after "General" ----
$names = @("Deployment", "Identity, Compliance and Config", "Maintain and Protect", "Applications")
foreach ($name in $names) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $ws.Name = $name
}

Write-Output "sheets added"
